$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G2").Value = 0.08788012380500777
$ws.Range("H2").Value = 32.51355987005816
$ws.Range("I2").Value = 12.9666078406142
$ws.Range("G3").Value = 0.1698029625776104
$ws.Range("H3").Value = 43.5809425366606
$ws.Range("G4").Value = -0.3140652136002539
$ws.Range("H4").Value = -14.20746299150609
$ws.Range("G5").Value = -0.4116723672468334
$ws.Range("H5").Value = -3.172188440973282
$ws.Range("G6").Value = 0.1934721856122315
$ws.Range("H6").Value = -1.865083048711867
$ws.Range("G7").Value = 0.2476447775230506
$ws.Range("H7").Value = 19.41462787740396
$ws.Range("G8").Value = 0.108905732269689
$ws.Range("H8").Value = 6.879282507407801
$ws.Range("G9").Value = 0.1603927918420133
$ws.Range("H9").Value = 26.81270887402953
$ws.Range("G10").Value = 0.07376548965297174
$ws.Range("H10").Value = 20.06530828788271
$ws.Range("G11").Value = 0.04866293974294369
$ws.Range("H11").Value = -2.537968190845442
$ws.Range("G12").Value = 0.09625283136901765
$ws.Range("H12").Value = 3.981561134935238
$ws.Range("G13").Value = 0.08314807663501843
$ws.Range("H13").Value = 9.107158821785363
$ws.Range("G14").Value = 0.1934106069158548
$ws.Range("H14").Value = -14.40555828037429
$ws.Range("G15").Value = 0.230177361217362
$ws.Range("H15").Value = -6.569668467593867
$ws.Range("G16").Value = 0.09854325898070752
$ws.Range("H16").Value = -13.36464720518527
$ws.Range("G17").Value = 0.124861488119469
$ws.Range("H17").Value = -16.42939707545014
$ws.Range("G18").Value = -0.009381490673663104
$ws.Range("H18").Value = -4.798301517552446
$ws.Range("G19").Value = -0.02703299766879391
$ws.Range("H19").Value = -211.6061139342761
$ws.Range("G20").Value = 0.152663515104255
$ws.Range("H20").Value = 79.47499904351368
$ws.Range("G21").Value = 0.1072516471118899
$ws.Range("H21").Value = 63.85813318919866
$ws.Range("G22").Value = 0.199689506267373
$ws.Range("H22").Value = 4.242929788442386
$ws.Range("G23").Value = 0.219941229812031
$ws.Range("H23").Value = 1.962801666580485
$ws.Range("G24").Value = 0.007493745544609639
$ws.Range("H24").Value = 296.9152646016782
$ws.Range("G25").Value = -0.00379867802964669
$ws.Range("H25").Value = 83.66674307620649
$ws.Range("G26").Value = 0.2188546325011767
$ws.Range("H26").Value = 6.828018902642608
$ws.Range("G27").Value = 0.1926788316990106
$ws.Range("H27").Value = -0.1066486457522578
$ws.Range("G28").Value = 0.03672226536399072
$ws.Range("H28").Value = -45.11993788692806
$ws.Range("G29").Value = 0.09232645547261413
$ws.Range("H29").Value = -2.056044956078971
